$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values
$ws.Range("B2").Value = 258.27178096320824
$ws.Range("C2").Value = 289.03401602611632
$ws.Range("D2").Value = 257.47125743493308
$ws.Range("E2").Value = 290.53438226129396

# Row 3 data values
$ws.Range("B3").Value = 259.59248883217685
$ws.Range("C3").Value = 287.87009579078614
$ws.Range("D3").Value = 254.3701899394438
$ws.Range("E3").Value = 297.45748213782667

# Update selection to match new range B1:E3
$ws.Range("B1:E3").Select()
